# Bugfixed QoQ Visualizations and a typo in the evaluation objects
#
# The first 16 data rows (rows 2-17, covering the earliest dates /
# 30864 .. 32234) are stale/incorrect and are removed from the series.
# Deleting these rows shifts every subsequent row up by 16, so the sheet
# ends up with 147 data rows (plus the header) instead of 163, i.e. the
# used range shrinks from A1:B164 to A1:B148.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A17").EntireRow.Delete()
